$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '275.60'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.186'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.581'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8260'

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01383'
$ws.Range("E10").Value = '9OneONE'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1627'
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08233'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03553'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03085'
$ws.Range("E14").Value = '13BitrueCoinBTR'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09127'
$ws.Range("E15").Value = '14BitMartTokenBMX'

# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.721'
$ws.Range("E16").Value = '15MCDexMCB'

# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001609'
$ws.Range("E17").Value = '16BitForexTokenBF'

# Row 18
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04700'
$ws.Range("E18").Value = '17CoinExTokenCET'

# Row 19
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006463'
$ws.Range("E19").Value = '18TigerCashTCH'

# Row 20
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001068'
$ws.Range("E20").Value = '19BitKanKAN'

# Row 21
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001503'
$ws.Range("E21").Value = '20NitroExNTX'

# Row 22
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.798'
$ws.Range("E22").Value = '21LEOLEO'

# Row 23
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.320'
$ws.Range("E23").Value = '22BTSETokenBTSE'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.006178'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0003748'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04677'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007016'

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004608'
$ws.Range("E42").Value = '41CEJICEJI'

# Row 43
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1101'
$ws.Range("E43").Value = '42BKEXTokenBKK'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01083'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006165'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.8464'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002602'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00001903'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01242'
